# Update "countries & provincias Spain" COVID stats sheet:
#  1) Refresh the numeric stats for the countries whose figures changed
#     (values are written by COUNTRY identity, at that country's
#     pre-refresh row position).
#  2) Re-sort the whole country table (A4:H219) by "Casos totales" (col B)
#     descending, matching the always-sorted presentation of this report
#     (a couple of countries swap places purely because their neighbour's
#     totals moved past them).
#  3) Bump the "Datos actualizados a ..." timestamp cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) refresh changed figures (row numbers are PRE-sort positions) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 4101308
$ws.Range("C4").Value = 433
$ws.Range("D4").Value = 1943503
$ws.Range("E4").Value = 2011613
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 146192

# Rusia (row 7)
$ws.Range("B7").Value = 795038
$ws.Range("C7").Value = 5848
$ws.Range("D7").Value = 580330
$ws.Range("E7").Value = 201816
$ws.Range("G7").Value = 147
$ws.Range("H7").Value = 12892

# Banglades (row 20)
$ws.Range("B20").Value = 216110
$ws.Range("C20").Value = 2856
$ws.Range("D20").Value = 119208
$ws.Range("E20").Value = 94101
$ws.Range("G20").Value = 50
$ws.Range("H20").Value = 2801

# Indonesia (row 27)
$ws.Range("B27").Value = 93657
$ws.Range("C27").Value = 1906
$ws.Range("D27").Value = 52164
$ws.Range("E27").Value = 36917
$ws.Range("G27").Value = 117
$ws.Range("H27").Value = 4576

# Kuwait (row 38) - figures unchanged, stays at its old value
$ws.Range("B38").Value = 61185
$ws.Range("D38").Value = 51520
$ws.Range("E38").Value = 9248
$ws.Range("H38").Value = 417

# Ucrania (row 39) - real update, overtakes Kuwait after the sort below
$ws.Range("B39").Value = 61851
$ws.Range("C39").Value = 856
$ws.Range("D39").Value = 34000
$ws.Range("E39").Value = 26300
$ws.Range("G39").Value = 17
$ws.Range("H39").Value = 1551

# Israel (row 41)
$ws.Range("B41").Value = 56748
$ws.Range("C41").Value = 663
$ws.Range("D41").Value = 23560
$ws.Range("E41").Value = 32755
$ws.Range("G41").Value = 3
$ws.Range("H41").Value = 433

# Polonia (row 47)
$ws.Range("B47").Value = 41580
$ws.Range("C47").Value = 418
$ws.Range("D47").Value = 31541
$ws.Range("E47").Value = 8388
$ws.Range("G47").Value = 9
$ws.Range("H47").Value = 1651

# Moldavia (row 63)
$ws.Range("D63").Value = 16174
$ws.Range("E63").Value = 4911
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 713

# Sri Lanka (row 114)
$ws.Range("D114").Value = 2077
$ws.Range("E114").Value = 664

# Eslovaquia (row 122)
$ws.Range("B122").Value = 2089
$ws.Range("C122").Value = 31
$ws.Range("E122").Value = 505

# Estonia (row 124) - small update, stays just behind Eslovenia after sort
$ws.Range("B124").Value = 2027
$ws.Range("C124").Value = 2
$ws.Range("D124").Value = 1912
$ws.Range("E124").Value = 46
$ws.Range("H124").Value = 69

# Eslovenia (row 125) - real update, overtakes Estonia after the sort below
$ws.Range("B125").Value = 2033
$ws.Range("C125").Value = 27
$ws.Range("D125").Value = 1648
$ws.Range("E125").Value = 270
$ws.Range("H125").Value = 115

# Letonia (row 138)
$ws.Range("B138").Value = 1203
$ws.Range("C138").Value = 6
$ws.Range("E138").Value = 127

# --- 2) re-sort the country table by "Casos totales" (col B) descending ---
$dataRange = $ws.Range("A4:H219")
$keyRange = $ws.Range("B4:B219")
$dataRange.Sort($keyRange, 2)

# --- 3) bump the "updated at" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 10:48"
